# Updated symbol list on Sun Dec 18 11:21:26 UTC 2022 with GitHub Actions
#
# Applies the per-cell value updates described by the upstream XML diff.
# Column D holds numeric-looking strings (prices) that must stay TEXT
# (leading/trailing zeros matter, e.g. "0.05620", "1.050"). Plain
# `.Value = "..."` would get auto-coerced to a number and lose that
# formatting, so for those cells we briefly force a text number format,
# assign the value, then restore the "Normal" style so no stray
# formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

function Set-TextualValue($cellRef, $val) {
    # Non-numeric-looking text (names, URLs, labels) - plain assignment
    # is safe and doesn't need the text-format trick.
    $ws.Range($cellRef).Value = $val
}

# --- Row 2 (BNB) ---
Set-TextValue "D2" "247.12"

# --- Row 4 (HuobiToken) ---
Set-TextValue "D4" "5.519"

# --- Row 5 (Cronos) ---
Set-TextValue "D5" "0.05620"

# --- Row 6 ---
Set-TextValue "D6" "3.374"

# --- Row 7 ---
Set-TextValue "D7" "6.473"

# --- Row 8 ---
Set-TextValue "D8" "0.8044"

# --- Row 9 ---
Set-TextValue "D9" "1.050"

# --- Row 11 ---
Set-TextValue "D11" "0.07258"

# --- Row 13 ---
Set-TextValue "D13" "0.02957"

# --- Row 14 ---
Set-TextValue "D14" "0.09259"

# --- Row 15 ---
Set-TextValue "D15" "0.001661"

# --- Row 16 ---
Set-TextValue "D16" "3.199"

# --- Row 17 ---
Set-TextValue "D17" "0.04693"

# --- Row 18 (One / ONE) ---
Set-TextValue "D18" "0.0005980"
Set-TextualValue "E18" "17OneONE"

# --- Row 19 ---
Set-TextValue "D19" "0.006270"

# --- Row 20 ---
Set-TextValue "D20" "0.001056"

# --- Row 21 ---
Set-TextValue "D21" "0.003817"

# --- Row 23 ---
Set-TextValue "D23" "0.0003305"

# --- Row 24 ---
Set-TextValue "D24" "3.972"

# --- Row 25 ---
Set-TextValue "D25" "2.130"

# --- Row 27 (ProBitToken / PROB) ---
Set-TextualValue "E27" "26ProBitTokenPROBBestin24h"

# --- Row 40 (IDEX) ---
Set-TextValue "D40" "0.04182"

# --- Row 41: was KickToken, now BKEXToken ---
Set-TextualValue "B41" "BKEXToken"
Set-TextualValue "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1044"
Set-TextualValue "E41" "40BKEXTokenBKK"

# --- Row 42 (CEJI) ---
Set-TextValue "D42" "0.002977"
Set-TextualValue "E42" "41CEJICEJI"

# --- Row 43: was BKEXToken, now KickToken ---
Set-TextualValue "B43" "KickToken"
Set-TextualValue "C43" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003249"
Set-TextualValue "E43" "42KickTokenKICKWorstin24h"

# --- Row 44 (LocalTraders) ---
Set-TextValue "D44" "0.009070"

# --- Row 45 (CoinLion) ---
Set-TextValue "D45" "0.00005646"

# --- Row 47 ---
Set-TextValue "D47" "0.6811"

# --- Row 48 ---
Set-TextValue "D48" "0.02568"

# --- Row 49 ---
Set-TextValue "D49" "0.00002103"
